$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("news")

# --- Insert new row 13 (pushes old row 13 "NTT法..." down to row 14) ---
$ws.Rows("13:13").Insert()

# --- Row 1: new header/info row above the table ---
$ws.Range("B1").Value = "議案審議情報"
$ws.Range("C1").Value = "https://www.sangiin.go.jp/japanese/joho1/kousei/gian/213/meisai/m213080213033.htm"
$ws.Rows("1:1").RowHeight = 39

# --- Row 13: new data row (Cabinet decision news item, dated 2024-03-01) ---
$ws.Range("A13").Value = (Get-Date -Year 2024 -Month 3 -Day 1)
$ws.Range("A13").NumberFormat = "m/d/yyyy"
$ws.Range("B13").Value = "令和6年3月1日（金）定例閣議案件"
$ws.Range("C13").Value = "https://www.kantei.go.jp/jp/kakugi/2024/kakugi-2024030101.html"
$ws.Rows("13:13").RowHeight = 30

# --- Hyperlinks for the two new URL cells ---
$ws.Hyperlinks.Add($ws.Range("C1"), "https://www.sangiin.go.jp/japanese/joho1/kousei/gian/213/meisai/m213080213033.htm")
$ws.Hyperlinks.Add($ws.Range("C13"), "https://www.kantei.go.jp/jp/kakugi/2024/kakugi-2024030101.html")

$ws.Range("A1:D1,A3:D14").Select()
